$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 used to hold the "Total Duration:" summary; it becomes a new
# punch-in/out data row, and the summary moves down to row 29.

# Write the new data row's values first. The date- and time-looking
# values are prefixed with an apostrophe so they are stored as literal
# text (matching the rest of the sheet) instead of being parsed into
# date/time serial numbers.
$ws.Range("A28").Value = "'2026-02-16"
$ws.Range("B28").Value = "'22:22:03"
$ws.Range("C28").Value = "23:22:45"
$ws.Range("D28").Value = "1.01 Hours"

# Copy the formatting (not the values) from the row above so row 28
# picks up the same cell style used throughout the table.
$ws.Range("A27:D27").Copy()
$ws.Range("A28:D28").PasteSpecial(-4122)

# Write the "Total Duration" summary into the new row 29.
$ws.Range("C29").Value = "Total Duration:"
$ws.Range("D29").Value = "33 Hours"

# Give row 29's summary cells the same formatting as the old summary
# row (now row 28's C:D cells), again copying formats only.
$ws.Range("C28:D28").Copy()
$ws.Range("C29:D29").PasteSpecial(-4122)

$excel.CutCopyMode = 0
